# "Add files via upload" — re-saved workbook that marks a handful of rows
# on the "Customers" and "Order Details" sheets as TRUE in a new
# "Test Result" (column F) flag, and leaves the "Order Details" sheet as
# the active/selected tab (previously "Products" was selected).

$wb = $excel.ActiveWorkbook

# --- Customers sheet: flag rows 2-4 in column F (Test Result) as TRUE ---
$wsCustomers = $wb.Worksheets.Item("Customers")
$wsCustomers.Range("F2:F4").Value = $true

# --- Order Details sheet: flag rows 2-10 in column F (Test Result) as TRUE ---
$wsOrderDetails = $wb.Worksheets.Item("Order Details")
$wsOrderDetails.Range("F2:F10").Value = $true

# Make "Order Details" the active sheet (tabSelected moves off "Products").
$wsOrderDetails.Activate()

Write-Output "Updated Test Result flags and activated Order Details sheet"
